$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal TEXT (not an auto-converted number),
# since these cells hold numbers formatted as plain text like "17.00".
function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.NumberFormat = "General"
    $rng.Style = "Normal"
}

# ECO Actual (B) and ECO Balance (D) columns updated with new figures.
Set-TextValue "B2" "17.00"
Set-TextValue "D2" "17.00"

Set-TextValue "B3" "8.00"
Set-TextValue "D3" "8.00"

Set-TextValue "B4" "24.00"
Set-TextValue "D4" "24.00"

Set-TextValue "B5" "9.00"
Set-TextValue "D5" "9.00"

Set-TextValue "B6" "9.00"
Set-TextValue "D6" "9.00"

Set-TextValue "B7" "67.00"
Set-TextValue "D7" "67.00"
